$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Status": fix casing of the "Seems fine..." status note, mark
# "Class 011 methods" as having the same note (previously "NONE"), and add a
# new note + highlight colour to "Cursor (RAMDAC_CU)".
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Status")

$newNote = "Seems fine, ALWAYS CHECK FOR BUGS"

# Simple text fixups (style/fill stay the same for these rows).
$wsStatus.Range("B5").Value  = $newNote
$wsStatus.Range("B10").Value = $newNote
$wsStatus.Range("B26").Value = $newNote

# Class 011 methods used to be "NONE" (red fill) -- now "Seems fine..." using
# the same green fill/style already used by B5/B10/B26. Copy that formatting
# across first so the existing style (fontId/fillId) gets reused instead of
# a brand new one being created, then overwrite the text.
$wsStatus.Range("B5").Copy($wsStatus.Range("B19")) | Out-Null
$wsStatus.Range("B19").Value = $newNote

# Cursor (RAMDAC_CU) previously had no status note -- add one and highlight
# it with a light-blue fill.
$wsStatus.Range("B34").Value = "There is special hardware support, but it just renders anyway using Class 011? At least emulate starting pos"
$wsStatus.Range("B34").Interior.Color = 15773696

# Column B needs to grow to fit the long new text.
$wsStatus.Columns("B").AutoFit() | Out-Null

# Leave column C selected on this sheet (as in the source file).
$wsStatus.Columns("C").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Bugs": two new known-broken items discovered while implementing
# patterns / bus work.
# ---------------------------------------------------------------------------
$wsBugs = $wb.Worksheets.Item("Bugs")

$wsBugs.Range("A12").Value = "AGP broken"
$wsBugs.Range("A13").Value = "ZX broken"

# Keep "Bugs" the active sheet/tab with A13 selected, matching the source.
$wsBugs.Range("A13").Select() | Out-Null
